$d = $word.ActiveDocument
$apost = [char]0x2019

# ------------------------------------------------------------------
# Locate the two paragraphs that matter:
#  - the intro paragraph "Almost every markdown application..." which
#    gets reworded and gains the "_GoBack" bookmark right after the
#    hyperlink to the spec
#  - the paragraph right under the "Application Support of Extended
#    Syntax" heading which currently owns that bookmark and needs to
#    lose it
# ------------------------------------------------------------------
$introIndex = -1
$bookmarkIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($introIndex -eq -1 -and $p.Range.Text -like "Almost every markdown application supports this syntax*") {
        $introIndex = $i
    }
    if ($bookmarkIndex -eq -1 -and $p.Range.WordOpenXML -like "*_GoBack*") {
        $bookmarkIndex = $i
    }
}

# ------------------------------------------------------------------
# 1. Drop the stale "_GoBack" bookmark from its old paragraph so the
#    id/name is free again.
# ------------------------------------------------------------------
if ($bookmarkIndex -ne -1) {
    $bp = $d.Paragraphs($bookmarkIndex)
    $bp.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>')
}

# ------------------------------------------------------------------
# 2. Rewrite the intro paragraph's runs (plain text only - no named
#    styles yet, those get re-applied below once the run boundaries
#    exist) and drop the "_GoBack" bookmark right after the hyperlink.
# ------------------------------------------------------------------
$introXml = (
  '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships">' +
  '<w:r><w:t xml:space="preserve">Almost every markdown application supports this syntax, which was </w:t></w:r>' +
  '<w:r><w:t xml:space="preserve">outlined in </w:t></w:r>' +
  '<w:r><w:t>John</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve"> Gruber' + $apost + 's</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
  '<w:hyperlink r:id="rId11" w:history="1">' +
    '<w:r><w:t xml:space="preserve">original </w:t></w:r>' +
    '<w:r><w:t>spec</w:t></w:r>' +
  '</w:hyperlink>' +
  '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
  '<w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">. </w:t></w:r>' +
  '<w:r><w:t>For a detailed</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve"> explanation, visit</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
  '<w:hyperlink r:id="rId12" w:history="1">' +
    '<w:r><w:t>here</w:t></w:r>' +
  '</w:hyperlink>' +
  '<w:r><w:t>.</w:t></w:r>' +
  '</w:p>'
)

$ip = $d.Paragraphs($introIndex)
$ip.Range.InsertXML($introXml)

# ------------------------------------------------------------------
# 3. Re-apply the "Hyperlink" character style to the runs that need
#    it: the trailing space right before the real hyperlink (which
#    also needs its underline switched off), and the two runs that
#    sit inside the hyperlink itself.
# ------------------------------------------------------------------
$ip = $d.Paragraphs($introIndex)
$sr = $ip.Range
$needle = "Gruber" + $apost + "s "
$null = $sr.Find.Execute($needle, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$sr.MoveStart(1, ($needle.Length - 1))
$sr.Style = "Hyperlink"
$sr.Underline = 0

$ip = $d.Paragraphs($introIndex)
$sr = $ip.Range
$null = $sr.Find.Execute("original ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$sr.Style = "Hyperlink"

$ip = $d.Paragraphs($introIndex)
$sr = $ip.Range
$null = $sr.Find.Execute("spec", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$sr.Style = "Hyperlink"

$ip = $d.Paragraphs($introIndex)
$sr = $ip.Range
$null = $sr.Find.Execute("here", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$sr.Style = "Hyperlink"
